$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# A1 Email, B1 Password, C1 Name stay the same.
# Reorder/replace D1:F1 and add a new G1 "ZipCode" header.
$ws.Range("D1").Value = "Phone"
$ws.Range("E1").Value = "DealerCode"
$ws.Range("F1").Value = "Address"
$ws.Range("G1").Value = "ZipCode"

# --- Data row (row 2) ---
# A2 (Email value) and C2 (Name value) stay the same.
# Password value changes from "Sai@123" to "Abc@123" and becomes a hyperlink.
# DealerCode value changes from "Test111" to "Test123".
# Add new ZipCode numeric value in G2.
$ws.Range("D2").Value = "(276) 343-7888"
$ws.Range("E2").Value = "Test123"
$ws.Range("F2").Value = "chandanagar,pjr statduim,xyz,hyderabad"
$ws.Range("G2").Value = 45765

# Turn B2 into a hyperlink showing "Abc@123" (replaces the old "Sai@123" text).
$ws.Range("B2").Value = "Abc@123"
$ws.Hyperlinks.Add($ws.Range("B2"), "Abc@123")

# Update the selected cell to match the authored selection.
$ws.Range("B2").Select() | Out-Null
